$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header style (copy bold/border/centered style from K1) to the new header cells L1:O1 ---
$ws.Range("K1").Copy()
$ws.Range("L1:O1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# --- Update header row text (G1:J1 keep their place but get renamed; K1 becomes "최종점수";
#     L1/M1/N1/O1 are new headers) ---
$ws.Range("G1").Value = "점수(룰)"
$ws.Range("H1").Value = "3일상승확률(%)"
$ws.Range("I1").Value = "5일상승확률(%)"
$ws.Range("J1").Value = "10일상승확률(%)"
$ws.Range("K1").Value = "최종점수"
$ws.Range("L1").Value = "예측방식"
$ws.Range("M1").Value = "판단"
$ws.Range("N1").Value = "MACRO_SCORE"
$ws.Range("O1").Value = "MACRO_SIGNAL"

# --- Column A holds date-looking text ("2025-11-29") that must stay plain text, not become
#     an Excel date serial number. Force text format first, write the values, then restore
#     the default (unstyled) look so the cells match a plain un-styled text cell. ---
$ws.Range("A2:A4").NumberFormat = "@"
$ws.Range("A2").Value = "2025-11-29"
$ws.Range("A3").Value = "2025-11-29"
$ws.Range("A4").Value = "2025-11-29"
$ws.Range("A2:A4").Style = "Normal"

# --- Row 2: Newmont Corporation / NEM ---
$ws.Range("B2").Value = "Newmont Corporation"
$ws.Range("C2").Value = "NEM"
$ws.Range("D2").Value = 90.79000000000001
$ws.Range("E2").ClearContents()
$ws.Range("F2").Value = 11.04
$ws.Range("G2").Value = 50
$ws.Range("H2").Value = 66
$ws.Range("I2").Value = 73
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 69.8
$ws.Range("L2").Value = "Pattern"
$ws.Range("M2").Value = "📈 매수 관찰 구간입니다."
$ws.Range("N2").Value = 85.36763896678245
$ws.Range("O2").Value = "🟢 완화적 (상승 우위)"

# --- Row 3: StreetTRACKS Gold Shares / GLD ---
$ws.Range("B3").Value = "StreetTRACKS Gold Shares"
$ws.Range("C3").Value = "GLD"
$ws.Range("D3").Value = 387.89
$ws.Range("E3").ClearContents()
$ws.Range("F3").Value = 3.48
$ws.Range("G3").Value = 50
$ws.Range("H3").Value = 56
$ws.Range("I3").Value = 70
$ws.Range("J3").Value = 73
$ws.Range("K3").Value = 68.6
$ws.Range("L3").Value = "Pattern"
$ws.Range("M3").Value = "📈 매수 관찰 구간입니다."
$ws.Range("N3").Value = 85.36763896678245
$ws.Range("O3").Value = "🟢 완화적 (상승 우위)"

# --- Row 4: Gold Dec 25 / GC=F ---
$ws.Range("B4").Value = "Gold Dec 25"
$ws.Range("C4").Value = "GC=F"
$ws.Range("D4").Value = 4247.7
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = 4.59
$ws.Range("G4").Value = 40
$ws.Range("H4").Value = 40
$ws.Range("I4").Value = 63
$ws.Range("J4").Value = 70
$ws.Range("K4").Value = 62.8
$ws.Range("L4").Value = "Pattern"
$ws.Range("M4").Value = "📈 매수 관찰 구간입니다."
$ws.Range("N4").Value = 85.36763896678245
$ws.Range("O4").Value = "🟢 완화적 (상승 우위)"
